# Cambio di segno alle derivate di controllo del latero direzionale
# Flip the sign of the lateral-directional control derivatives
# (columns E = delta_r, F = delta_a) for rows 10-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 10..15) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$r")
        $cell.Value = -1 * $cell.Value2
    }
}
